# Update gh-pages generated output numbers (views/stats counters)
# Sheet "展览" (Exhibition): F3, F5, F7, F9
# Sheet "全部类型" (All types): F3, F5, F9, F11

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 197
$wsExhibition.Range("F5").Value = 452
$wsExhibition.Range("F7").Value = 2478
$wsExhibition.Range("F9").Value = 6532

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 197
$wsAll.Range("F5").Value = 452
$wsAll.Range("F9").Value = 2478
$wsAll.Range("F11").Value = 6532
